$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$genesToDelete = @("Adrb1", "Adrb2", "Adrb3", "Lipe", "Lpl", "Pde3b")

for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    $cellText = $row.Cells.Item(1).Range.Text
    foreach ($gene in $genesToDelete) {
        if ($cellText -like "$gene*") {
            $row.Delete()
            break
        }
    }
}
